# Revert "Migration to Automation-Org/TestCases-maintenance/WIP-RMA TestCases"
#
# 1) Remove the "Receipt" worksheet that the migration had added.
# 2) Restore the RMA-VS7X-* test identifiers on the "RMA Details Maintenance
#    Grid" sheet (the migration had overwritten them with RMA-UTSM-* values).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$grid = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$grid.Range("E2").Value = "RMA-VS7X-002"
$grid.Range("F2").Value = "RMA-VS7X-1-2"
$grid.Range("J2").Value = "a6h1K000000PrnUQAS"

$grid.Range("E3").Value = "RMA-VS7X-001"
$grid.Range("F3").Value = "RMA-VS7X-1-1"
$grid.Range("J3").Value = "a6h1K000000PrnTQAS"

$grid.Range("E4").Value = "RMA-VS7X-003"
$grid.Range("F4").Value = "RMA-VS7X-1-3"
$grid.Range("J4").Value = "a6h1K000000PrnVQAS"

$receipt = $wb.Worksheets.Item("Receipt")
$receipt.Delete()
